# Adding update, interests and article
$wb = $excel.ActiveWorkbook

# --- Organizations sheet: add two new organizations with hyperlinks ---
$orgs = $wb.Worksheets.Item("Organizations")

$orgs.Cells.Item(18, 1).Value = 17
$orgs.Cells.Item(18, 2).Value = "ideas42"
$orgs.Hyperlinks.Add($orgs.Cells.Item(18, 2), "http://www.ideas42.org/") | Out-Null
$orgs.Cells.Item(18, 2).Style = "Hyperlink"

$orgs.Cells.Item(19, 1).Value = 18
$orgs.Cells.Item(19, 2).Value = "Regulation Room"
$orgs.Hyperlinks.Add($orgs.Cells.Item(19, 2), "http://regulationroom.org/") | Out-Null
$orgs.Cells.Item(19, 2).Style = "Hyperlink"

$orgs.Range("A20").Select()

# --- Make Organizations the active/selected sheet and tab ---
$orgs.Activate()

$wb.Save()
